# Adicionados balancos concatenados em uma unica planilha.
# Extends VITT3 sheet with three additional quarterly columns: R (31/12/2023), S (31/03/2024), T (30/06/2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from Q1 onto the new header cells R1:T1
$ws.Range("Q1").Copy()
$ws.Range("R1:T1").PasteSpecial(-4122)

# New quarter headers
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

$ws.Range("R2").Value = 925756.992
$ws.Range("S2").Value = 917881.9840000001
$ws.Range("T2").Value = 847132.992

$ws.Range("R3").Value = 571633.024
$ws.Range("S3").Value = 563324.992
$ws.Range("T3").Value = 476673.984

$ws.Range("R4").Value = 82829
$ws.Range("S4").Value = 17640
$ws.Range("T4").Value = 38793

$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 70
$ws.Range("T5").Value = 2301

$ws.Range("R6").Value = 304176.992
$ws.Range("S6").Value = 346944.992
$ws.Range("T6").Value = 182066

$ws.Range("R7").Value = 154336.992
$ws.Range("S7").Value = 168791.008
$ws.Range("T7").Value = 212294

$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0

$ws.Range("R9").Value = 11496
$ws.Range("S9").Value = 9828
$ws.Range("T9").Value = 12519

$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0

$ws.Range("R11").Value = 18794
$ws.Range("S11").Value = 20051
$ws.Range("T11").Value = 28701

$ws.Range("R12").Value = 27480
$ws.Range("S12").Value = 25922
$ws.Range("T12").Value = 37124

$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0

$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0

$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0

$ws.Range("R16").Value = 1808
$ws.Range("S16").Value = 3065
$ws.Range("T16").Value = 7634

$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0

$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 0

$ws.Range("R19").Value = 15630
$ws.Range("S19").Value = 16306
$ws.Range("T19").Value = 25172

$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0

$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0

$ws.Range("R22").Value = 256
$ws.Range("S22").Value = 253
$ws.Range("T22").Value = 255

$ws.Range("R23").Value = 308955.008
$ws.Range("S23").Value = 311688.992
$ws.Range("T23").Value = 316788.992

$ws.Range("R24").Value = 17433
$ws.Range("S24").Value = 16693
$ws.Range("T24").Value = 16291

$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 0

$ws.Range("R26").Value = 925756.992
$ws.Range("S26").Value = 917881.9840000001
$ws.Range("T26").Value = 847132.992

$ws.Range("R27").Value = 277656.992
$ws.Range("S27").Value = 285696
$ws.Range("T27").Value = 194514

$ws.Range("R28").Value = 17908
$ws.Range("S28").Value = 19852
$ws.Range("T28").Value = 25874

$ws.Range("R29").Value = 16734
$ws.Range("S29").Value = 34882
$ws.Range("T29").Value = 28347

$ws.Range("R30").Value = 4434
$ws.Range("S30").Value = 1445
$ws.Range("T30").Value = 1841

$ws.Range("R31").Value = 172147.008
$ws.Range("S31").Value = 195492.992
$ws.Range("T31").Value = 108903

$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 0
$ws.Range("T32").Value = 0

$ws.Range("R33").Value = 38297
$ws.Range("S33").Value = 15809
$ws.Range("T33").Value = 106

$ws.Range("R34").Value = 28137
$ws.Range("S34").Value = 18215
$ws.Range("T34").Value = 29443

$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 0
$ws.Range("T35").Value = 0

$ws.Range("R36").Value = 0
$ws.Range("S36").Value = 0
$ws.Range("T36").Value = 0

$ws.Range("R37").Value = 27728
$ws.Range("S37").Value = 26999
$ws.Range("T37").Value = 80725

$ws.Range("R38").Value = 27391
$ws.Range("S38").Value = 26731
$ws.Range("T38").Value = 80138

$ws.Range("R39").Value = 0
$ws.Range("S39").Value = 0
$ws.Range("T39").Value = 0

$ws.Range("R40").Value = 289
$ws.Range("S40").Value = 0
$ws.Range("T40").Value = 0

$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 0
$ws.Range("T41").Value = 0

$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("T42").Value = 0

$ws.Range("R43").Value = 48
$ws.Range("S43").Value = 268
$ws.Range("T43").Value = 587

$ws.Range("R44").Value = 0
$ws.Range("S44").Value = 0
$ws.Range("T44").Value = 0

$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("T45").Value = 0

$ws.Range("R46").Value = 2
$ws.Range("S46").Value = 94
$ws.Range("T46").Value = 69

$ws.Range("R47").Value = 620369.968
$ws.Range("S47").Value = 605093.008
$ws.Range("T47").Value = 571825.0159999999

$ws.Range("R48").Value = 255208.992
$ws.Range("S48").Value = 465640.992
$ws.Range("T48").Value = 465640.992

$ws.Range("R49").Value = -10203
$ws.Range("S49").Value = -5381
$ws.Range("T49").Value = -21392

$ws.Range("R50").Value = 1768
$ws.Range("S50").Value = 1711
$ws.Range("T50").Value = 2129

$ws.Range("R51").Value = 373596
$ws.Range("S51").Value = 142311.008
$ws.Range("T51").Value = 142311.008

$ws.Range("R52").Value = 0
$ws.Range("S52").Value = 811
$ws.Range("T52").Value = -16864

$ws.Range("R53").Value = 0
$ws.Range("S53").Value = 0
$ws.Range("T53").Value = 0

$ws.Range("R54").Value = 0
$ws.Range("S54").Value = 0
$ws.Range("T54").Value = 0

$ws.Range("R55").Value = 0
$ws.Range("S55").Value = 0
$ws.Range("T55").Value = 0

$ws.Range("R56").Value = 0
$ws.Range("S56").Value = 0
$ws.Range("T56").Value = 0

# Row 57 is a separator row with no data; R57:T57 stay blank
# Row 58 is a separator row with no data; R58:T58 stay blank
$ws.Range("R59").Value = 243298
$ws.Range("S59").Value = 121559
$ws.Range("T59").Value = 99893

$ws.Range("R60").Value = -157476.976
$ws.Range("S60").Value = -77523
$ws.Range("T60").Value = -82524

$ws.Range("R61").Value = 85821.008
$ws.Range("S61").Value = 44036
$ws.Range("T61").Value = 17369

$ws.Range("R62").Value = -26340
$ws.Range("S62").Value = -20669
$ws.Range("T62").Value = -17622

$ws.Range("R63").Value = -28012
$ws.Range("S63").Value = -25356
$ws.Range("T63").Value = -27129

$ws.Range("R64").Value = 1830
$ws.Range("S64").Value = 289
$ws.Range("T64").Value = 990

$ws.Range("R65").Value = -60
$ws.Range("S65").Value = -30
$ws.Range("T65").Value = 154

$ws.Range("R66").Value = 0
$ws.Range("S66").Value = 0
$ws.Range("T66").Value = 0

$ws.Range("R67").Value = 0
$ws.Range("S67").Value = 0
$ws.Range("T67").Value = 0

$ws.Range("R68").Value = -2343
$ws.Range("S68").Value = 2015
$ws.Range("T68").Value = 76

$ws.Range("R69").Value = 16785
$ws.Range("S69").Value = 8837
$ws.Range("T69").Value = 10900

$ws.Range("R70").Value = -19128
$ws.Range("S70").Value = -6822
$ws.Range("T70").Value = -10824

# Row 71 is a separator row with no data; R71:T71 stay blank
# Row 72 is a separator row with no data; R72:T72 stay blank
# Row 73 is a separator row with no data; R73:T73 stay blank
$ws.Range("R74").Value = 30896
$ws.Range("S74").Value = 285
$ws.Range("T74").Value = -26162

$ws.Range("R75").Value = 6881
$ws.Range("S75").Value = -161
$ws.Range("T75").Value = -413

$ws.Range("R76").Value = 3486
$ws.Range("S76").Value = 676
$ws.Range("T76").Value = 8867

# Row 77 is a separator row with no data; R77:T77 stay blank
# Row 78 is a separator row with no data; R78:T78 stay blank
$ws.Range("R79").Value = 52
$ws.Range("S79").Value = 11
$ws.Range("T79").Value = 33

$ws.Range("R80").Value = 41315
$ws.Range("S80").Value = 811
$ws.Range("T80").Value = -17675

